$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_9_4_4"
$ws.Range("B2").Value = 0.2112709723341609
$ws.Range("C2").Value = -1.270970585416494
$ws.Range("D2").Value = -1.480401527222878
$ws.Range("E2").Value = -1.349045258114328
$ws.Range("F2").Value = 0.8728905916213989
$ws.Range("G2").Value = 2.136618852615356
$ws.Range("H2").Value = 2.860383749008179
$ws.Range("I2").Value = 2.477213859558105
$ws.Range("A3").Value = "model_9_4_3"
$ws.Range("B3").Value = 0.2616871498882661
$ws.Range("C3").Value = -1.196313221554341
$ws.Range("D3").Value = -1.086603741683489
$ws.Range("E3").Value = -1.111133663315895
$ws.Range("F3").Value = 0.8170947432518005
$ws.Range("G3").Value = 2.066378355026245
$ws.Range("H3").Value = 2.406258583068848
$ws.Range("I3").Value = 2.226321220397949
$ws.Range("A4").Value = "model_9_4_12"
$ws.Range("B4").Value = 0.2646761391190556
$ws.Range("C4").Value = -2.382229510258318
$ws.Range("D4").Value = -0.5944832923852155
$ws.Range("E4").Value = -1.41802081217306
$ws.Range("F4").Value = 0.8137868046760559
$ws.Range("G4").Value = 3.182135343551636
$ws.Range("H4").Value = 1.838748216629028
$ws.Range("I4").Value = 2.549952745437622
$ws.Range("A5").Value = "model_9_4_13"
$ws.Range("B5").Value = 0.2774840200016251
$ws.Range("C5").Value = -2.257340757991858
$ws.Range("D5").Value = -0.6360018472065014
$ws.Range("E5").Value = -1.380401035811197
$ws.Range("F5").Value = 0.7996123433113098
$ws.Range("G5").Value = 3.064635038375854
$ws.Range("H5").Value = 1.886627197265625
$ws.Range("I5").Value = 2.510280609130859
$ws.Range("A6").Value = "model_9_4_14"
$ws.Range("B6").Value = 0.2867112208403876
$ws.Range("C6").Value = -2.169912727645614
$ws.Range("D6").Value = -0.665954933875883
$ws.Range("E6").Value = -1.354517506240311
$ws.Range("F6").Value = 0.7894005179405212
$ws.Range("G6").Value = 2.98237943649292
$ws.Range("H6").Value = 1.921168923377991
$ws.Range("I6").Value = 2.48298454284668
$ws.Range("A7").Value = "model_9_4_16"
$ws.Range("B7").Value = 0.3154815794925628
$ws.Range("C7").Value = -1.971187127413387
$ws.Range("D7").Value = -0.6739542430528378
$ws.Range("E7").Value = -1.264773335350963
$ws.Range("F7").Value = 0.7575601935386658
$ws.Range("G7").Value = 2.795410394668579
$ws.Range("H7").Value = 1.930393576622009
$ws.Range("I7").Value = 2.388343811035156
$ws.Range("A8").Value = "model_9_4_15"
$ws.Range("B8").Value = 0.3184152598807807
$ws.Range("C8").Value = -1.942190852363367
$ws.Range("D8").Value = -0.6740503705522891
$ws.Range("E8").Value = -1.251124988145825
$ws.Range("F8").Value = 0.7543134689331055
$ws.Range("G8").Value = 2.768129587173462
$ws.Range("H8").Value = 1.930504560470581
$ws.Range("I8").Value = 2.373950958251953
$ws.Range("A9").Value = "model_9_4_11"
$ws.Range("B9").Value = 0.3233735469341812
$ws.Range("C9").Value = -2.088618753453285
$ws.Range("D9").Value = -0.4507237700830002
$ws.Range("E9").Value = -1.205366007704
$ws.Range("F9").Value = 0.7488260269165039
$ws.Range("G9").Value = 2.90589451789856
$ws.Range("H9").Value = 1.6729656457901
$ws.Range("I9").Value = 2.325695037841797
$ws.Range("A10").Value = "model_9_4_21"
$ws.Range("B10").Value = 0.3346819351570769
$ws.Range("C10").Value = -1.90883458265419
$ws.Range("D10").Value = -0.6282244328089379
$ws.Range("E10").Value = -1.211792477566435
$ws.Range("F10").Value = 0.7363110780715942
$ws.Range("G10").Value = 2.736746549606323
$ws.Range("H10").Value = 1.877658367156982
$ws.Range("I10").Value = 2.332472324371338
$ws.Range("A11").Value = "model_9_4_20"
$ws.Range("B11").Value = 0.335445881839046
$ws.Range("C11").Value = -1.897589045771548
$ws.Range("D11").Value = -0.6300533993988133
$ws.Range("E11").Value = -1.207421499451616
$ws.Range("F11").Value = 0.7354655861854553
$ws.Range("G11").Value = 2.726166486740112
$ws.Range("H11").Value = 1.879767417907715
$ws.Range("I11").Value = 2.327862739562988
$ws.Range("A12").Value = "model_9_4_19"
$ws.Range("B12").Value = 0.3356357858558744
$ws.Range("C12").Value = -1.893001178798407
$ws.Range("D12").Value = -0.630754448310989
$ws.Range("E12").Value = -1.205614355625152
$ws.Range("F12").Value = 0.7352553606033325
$ws.Range("G12").Value = 2.721850156784058
$ws.Range("H12").Value = 1.880575895309448
$ws.Range("I12").Value = 2.325957059860229
$ws.Range("A13").Value = "model_9_4_18"
$ws.Range("B13").Value = 0.3361591780354206
$ws.Range("C13").Value = -1.864737089268072
$ws.Range("D13").Value = -0.6462473979726242
$ws.Range("E13").Value = -1.200232699902365
$ws.Range("F13").Value = 0.7346762418746948
$ws.Range("G13").Value = 2.695258140563965
$ws.Range("H13").Value = 1.898442387580872
$ws.Range("I13").Value = 2.320281982421875
$ws.Range("A14").Value = "model_9_4_22"
$ws.Range("B14").Value = 0.3382844936628217
$ws.Range("C14").Value = -1.875783251074439
$ws.Range("D14").Value = -0.6346403570650616
$ws.Range("E14").Value = -1.199484893069932
$ws.Range("F14").Value = 0.7323240637779236
$ws.Range("G14").Value = 2.705650806427002
$ws.Range("H14").Value = 1.885057210922241
$ws.Range("I14").Value = 2.319493055343628
$ws.Range("A15").Value = "model_9_4_17"
$ws.Range("B15").Value = 0.3436479123795986
$ws.Range("C15").Value = -1.789972321721864
$ws.Range("D15").Value = -0.6570393566774793
$ws.Range("E15").Value = -1.170477512040379
$ws.Range("F15").Value = 0.7263883948326111
$ws.Range("G15").Value = 2.624916315078735
$ws.Range("H15").Value = 1.910887479782104
$ws.Range("I15").Value = 2.28890323638916
$ws.Range("A16").Value = "model_9_4_23"
$ws.Range("B16").Value = 0.3464886290860391
$ws.Range("C16").Value = -1.850234332281612
$ws.Range("D16").Value = -0.6067562045531472
$ws.Range("E16").Value = -1.173066241008569
$ws.Range("F16").Value = 0.7232445478439331
$ws.Range("G16").Value = 2.681613206863403
$ws.Range("H16").Value = 1.852901339530945
$ws.Range("I16").Value = 2.291633129119873
$ws.Range("A17").Value = "model_9_4_24"
$ws.Range("B17").Value = 0.3514776630744696
$ws.Range("C17").Value = -1.82502107851598
$ws.Range("D17").Value = -0.5981059047655606
$ws.Range("E17").Value = -1.156703134050209
$ws.Range("F17").Value = 0.7177231907844543
$ws.Range("G17").Value = 2.657891511917114
$ws.Range("H17").Value = 1.842925906181335
$ws.Range("I17").Value = 2.274377346038818
$ws.Range("A18").Value = "model_9_4_2"
$ws.Range("B18").Value = 0.3646554409915991
$ws.Range("C18").Value = -0.9894190467782642
$ws.Range("D18").Value = -0.3859966234800596
$ws.Range("E18").Value = -0.6528798979742243
$ws.Range("F18").Value = 0.7031391859054565
$ws.Range("G18").Value = 1.871724009513855
$ws.Range("H18").Value = 1.598322749137878
$ws.Range("I18").Value = 1.743064403533936
$ws.Range("A19").Value = "model_9_4_10"
$ws.Range("B19").Value = 0.3654195420909396
$ws.Range("C19").Value = -1.919168839476508
$ws.Range("D19").Value = -0.2954449990472214
$ws.Range("E19").Value = -1.045422115872919
$ws.Range("F19").Value = 0.7022935748100281
$ws.Range("G19").Value = 2.746469497680664
$ws.Range("H19").Value = 1.49389922618866
$ws.Range("I19").Value = 2.157024621963501
$ws.Range("A20").Value = "model_9_4_9"
$ws.Range("B20").Value = 0.3727127498771917
$ws.Range("C20").Value = -1.875876514206552
$ws.Range("D20").Value = -0.2860823499271949
$ws.Range("E20").Value = -1.020156635141773
$ws.Range("F20").Value = 0.6942221522331238
$ws.Range("G20").Value = 2.705738306045532
$ws.Range("H20").Value = 1.483102202415466
$ws.Range("I20").Value = 2.130380392074585
$ws.Range("A21").Value = "model_9_4_1"
$ws.Range("B21").Value = 0.3761346146798321
$ws.Range("C21").Value = -1.100135983982577
$ws.Range("D21").Value = -0.004470559487531833
$ws.Range("E21").Value = -0.5088407156970844
$ws.Range("F21").Value = 0.6904352307319641
$ws.Range("G21").Value = 1.97589111328125
$ws.Range("H21").Value = 1.158349275588989
$ws.Range("I21").Value = 1.591166138648987
$ws.Range("A22").Value = "model_9_4_8"
$ws.Range("B22").Value = 0.3807274360988174
$ws.Range("C22").Value = -1.823351179332678
$ws.Range("D22").Value = -0.2754137738593732
$ws.Range("E22").Value = -0.9898585463828951
$ws.Range("F22").Value = 0.6853522658348083
$ws.Range("G22").Value = 2.656320333480835
$ws.Range("H22").Value = 1.470799326896667
$ws.Range("I22").Value = 2.098428964614868
$ws.Range("A23").Value = "model_9_4_0"
$ws.Range("B23").Value = 0.3845605440001083
$ws.Range("C23").Value = -1.044387467141385
$ws.Range("D23").Value = 0.06467207257053653
$ws.Range("E23").Value = -0.4469285963205309
$ws.Range("F23").Value = 0.6811100840568542
$ws.Range("G23").Value = 1.92344069480896
$ws.Range("H23").Value = 1.078614354133606
$ws.Range("I23").Value = 1.525876045227051
$ws.Range("A24").Value = "model_9_4_7"
$ws.Range("B24").Value = 0.4130598243683032
$ws.Range("C24").Value = -1.629184744699681
$ws.Range("D24").Value = -0.2243203398982363
$ws.Range("E24").Value = -0.8718558751648666
$ws.Range("F24").Value = 0.6495698094367981
$ws.Range("G24").Value = 2.473641157150269
$ws.Range("H24").Value = 1.41187858581543
$ws.Range("I24").Value = 1.973988056182861
$ws.Range("A25").Value = "model_9_4_6"
$ws.Range("B25").Value = 0.423366998467543
$ws.Range("C25").Value = -1.572081434631954
$ws.Range("D25").Value = -0.1651995271361948
$ws.Range("E25").Value = -0.8144596324301672
$ws.Range("F25").Value = 0.6381627917289734
$ws.Range("G25").Value = 2.419915914535522
$ws.Range("H25").Value = 1.343700885772705
$ws.Range("I25").Value = 1.913460254669189
$ws.Range("A26").Value = "model_9_4_5"
$ws.Range("B26").Value = 0.4236467407323959
$ws.Range("C26").Value = -1.568213221280754
$ws.Range("D26").Value = -0.1661300269616122
$ws.Range("E26").Value = -0.813112130595113
$ws.Range("F26").Value = 0.63785320520401
$ws.Range("G26").Value = 2.416276454925537
$ws.Range("H26").Value = 1.885057210922241
$ws.Range("I26").Value = 1.912039160728455
